# Scheduled-runner update: refresh market-derived profit figures (columns H-N)
# on the leve rows whose current Marketboard prices moved since the last sync.
# Values below were supplied by the scheduled data-refresh job.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Range("H17").Value = 6251966.5
$ws.Range("J17").Value = 6668637.5
$ws.Range("L17").Value = 20005912.5
$ws.Range("N17").Value = -20006248.5

# Row 92: Whinier than the Sword
$ws.Range("H92").Value = 567.9231
$ws.Range("I92").Value = 571.1818
$ws.Range("J92").Value = 550
$ws.Range("K92").Value = 571.1818
$ws.Range("L92").Value = 550
$ws.Range("M92").Value = 676.8182
$ws.Range("N92").Value = -3046

# Row 96: Scroll Down
$ws.Range("H96").Value = 19231444
$ws.Range("I96").Value = 27778408
$ws.Range("K96").Value = 83335224
$ws.Range("M96").Value = -83333851

# Row 106: Making Your Mark
$ws.Range("H106").Value = 10755838
$ws.Range("I106").Value = 41668430
$ws.Range("K106").Value = 41668430
$ws.Range("M106").Value = -41667799

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1060.1621
$ws.Range("I137").Value = 967.7143
$ws.Range("K137").Value = 2903.1429
$ws.Range("M137").Value = -353.1428999999998

# Row 138: All-night Crafting
$ws.Range("H138").Value = 2118.1755
$ws.Range("I138").Value = 1313.8096
$ws.Range("J138").Value = 2587.389
$ws.Range("K138").Value = 3941.4288
$ws.Range("L138").Value = 7762.167
$ws.Range("M138").Value = 1198.5712
$ws.Range("N138").Value = -18042.167

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 6471.1865
$ws.Range("I32").Value = 5248.864
$ws.Range("J32").Value = 15434.889
$ws.Range("K32").Value = 5248.864
$ws.Range("L32").Value = 15434.889
$ws.Range("M32").Value = -4961.864
$ws.Range("N32").Value = -16008.889

# Row 97: Ore for Me
$ws.Range("H97").Value = 111112750
$ws.Range("I97").Value = 1499.8334
$ws.Range("K97").Value = 1499.8334
$ws.Range("M97").Value = -1003.8334

# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 1926.6666
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 1890
$ws.Range("K102").Value = 2000
$ws.Range("L102").Value = 1890
$ws.Range("M102").Value = -378
$ws.Range("N102").Value = -5134

$ws = $wb.Worksheets.Item("BSM")
# Row 19: Twice as Slice
$ws.Range("H19").Value = 600
$ws.Range("I19").Value = 600
$ws.Range("K19").Value = 600
$ws.Range("M19").Value = -427

# Row 94: High Steal
$ws.Range("H94").Value = 696.0357
$ws.Range("I94").Value = 656.6
$ws.Range("J94").Value = 741.53845
$ws.Range("K94").Value = 656.6
$ws.Range("L94").Value = 741.53845
$ws.Range("M94").Value = -205.6
$ws.Range("N94").Value = -1643.53845

# Row 99: Meddle in Metal
$ws.Range("H99").Value = 1052.75
$ws.Range("I99").Value = 733.3333
$ws.Range("K99").Value = 733.3333
$ws.Range("M99").Value = 764.6667

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 2382824.2
$ws.Range("J105").Value = 2779645
$ws.Range("L105").Value = 2779645
$ws.Range("N105").Value = -2783139

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 3592.3953
$ws.Range("I134").Value = 3561.75
$ws.Range("K134").Value = 10685.25
$ws.Range("M134").Value = -8150.25

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 2985.5962
$ws.Range("I31").Value = 1559.7333
$ws.Range("J31").Value = 4929.9546
$ws.Range("K31").Value = 1559.7333
$ws.Range("L31").Value = 4929.9546
$ws.Range("M31").Value = -1264.7333
$ws.Range("N31").Value = -5519.9546

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2985.5962
$ws.Range("I34").Value = 1559.7333
$ws.Range("J34").Value = 4929.9546
$ws.Range("K34").Value = 1559.7333
$ws.Range("L34").Value = 4929.9546
$ws.Range("M34").Value = -1357.7333
$ws.Range("N34").Value = -5333.9546

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 37828.715
$ws.Range("I58").Value = 2482
$ws.Range("J58").Value = 73175.42999999999
$ws.Range("K58").Value = 2482
$ws.Range("L58").Value = 73175.42999999999
$ws.Range("M58").Value = -2279
$ws.Range("N58").Value = -73581.42999999999

# Row 99: O Pine
$ws.Range("H99").Value = 3194.423
$ws.Range("I99").Value = 2686.389
$ws.Range("K99").Value = 2686.389
$ws.Range("M99").Value = -1188.389

# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 1137.7693
$ws.Range("I105").Value = 989.0909
$ws.Range("K105").Value = 989.0909
$ws.Range("M105").Value = 757.9091

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 1718
$ws.Range("I122").Value = 1597.1428
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4791.428400000001
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2341.428400000001
$ws.Range("N122").Value = -10900

# Row 126: A Better Conductor
$ws.Range("H126").Value = 3194.423
$ws.Range("I126").Value = 2686.389
$ws.Range("K126").Value = 8059.167
$ws.Range("M126").Value = -5589.167

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 1533.5646
$ws.Range("I132").Value = 1204.6666
$ws.Range("J132").Value = 2661.2144
$ws.Range("K132").Value = 3613.9998
$ws.Range("L132").Value = 7983.6432
$ws.Range("M132").Value = -1083.9998
$ws.Range("N132").Value = -13043.6432

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 984.61536
$ws.Range("I134").Value = 917.2
$ws.Range("J134").Value = 1209.3334
$ws.Range("K134").Value = 2751.6
$ws.Range("L134").Value = 3628.0002
$ws.Range("M134").Value = -216.6000000000004
$ws.Range("N134").Value = -8698.0002

# Row 136: Turali Quality
$ws.Range("H136").Value = 37828.715
$ws.Range("I136").Value = 2482
$ws.Range("J136").Value = 73175.42999999999
$ws.Range("K136").Value = 7446
$ws.Range("L136").Value = 219526.29
$ws.Range("M136").Value = -4896
$ws.Range("N136").Value = -224626.29

$ws = $wb.Worksheets.Item("CUL")
# Row 60: Drinking to Your Health
$ws.Range("H60").Value = 1500
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 1500
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 4500
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -5002

# Row 97: The Frier Never Lies
$ws.Range("H97").Value = 711.8
$ws.Range("J97").Value = 711.8
$ws.Range("L97").Value = 2135.4
$ws.Range("N97").Value = -3127.4

# Row 105: Fish Box
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 682.4
$ws.Range("J131").Value = 709.2308
$ws.Range("L131").Value = 2127.6924
$ws.Range("N131").Value = -12207.6924

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 4177061
$ws.Range("J70").Value = 7827938.5
$ws.Range("L70").Value = 7827938.5
$ws.Range("N70").Value = -7828478.5

# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 4177061
$ws.Range("J73").Value = 7827938.5
$ws.Range("L73").Value = 7827938.5
$ws.Range("N73").Value = -7829810.5

# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 1711
$ws.Range("I97").Value = 1755.4166
$ws.Range("J97").Value = 1533.3334
$ws.Range("K97").Value = 1755.4166
$ws.Range("L97").Value = 1533.3334
$ws.Range("M97").Value = -1259.4166
$ws.Range("N97").Value = -2525.3334

# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 2174.9092
$ws.Range("I102").Value = 1841.2632
$ws.Range("J102").Value = 4288
$ws.Range("K102").Value = 1841.2632
$ws.Range("L102").Value = 4288
$ws.Range("M102").Value = -219.2632000000001
$ws.Range("N102").Value = -7532

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 2946.868
$ws.Range("I126").Value = 2457.9355
$ws.Range("J126").Value = 3635.818
$ws.Range("K126").Value = 7373.806500000001
$ws.Range("L126").Value = 10907.454
$ws.Range("M126").Value = -4903.806500000001
$ws.Range("N126").Value = -15847.454

# Row 132: On Board for Lar
$ws.Range("H132").Value = 19591.129
$ws.Range("I132").Value = 3476.72
$ws.Range("K132").Value = 10430.16
$ws.Range("M132").Value = -7900.16

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad
$ws.Range("H40").Value = 2803.6333
$ws.Range("I40").Value = 2531.4783
$ws.Range("J40").Value = 3697.8572
$ws.Range("K40").Value = 2531.4783
$ws.Range("L40").Value = 3697.8572
$ws.Range("M40").Value = -2395.4783
$ws.Range("N40").Value = -3969.8572

# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 7222.5557
$ws.Range("I61").Value = 2600.8
$ws.Range("J61").Value = 12999.75
$ws.Range("K61").Value = 2600.8
$ws.Range("L61").Value = 12999.75
$ws.Range("M61").Value = -2398.8
$ws.Range("N61").Value = -13403.75

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 1595.8572
$ws.Range("I93").Value = 1595.8572
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1595.8572
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -347.8571999999999
$ws.Range("N93").ClearContents()

# Row 100: Tiger in the Sack
$ws.Range("H100").Value = 2327.1428
$ws.Range("I100").Value = 1763.6666
$ws.Range("J100").Value = 2749.75
$ws.Range("K100").Value = 1763.6666
$ws.Range("L100").Value = 2749.75
$ws.Range("M100").Value = -1222.6666
$ws.Range("N100").Value = -3831.75

# Row 113: Peace in Rest
$ws.Range("H113").Value = 7222.5557
$ws.Range("I113").Value = 2600.8
$ws.Range("J113").Value = 12999.75
$ws.Range("K113").Value = 2600.8
$ws.Range("L113").Value = 12999.75
$ws.Range("M113").Value = -430.8000000000002
$ws.Range("N113").Value = -17339.75

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 336420.88
$ws.Range("I132").Value = 524969.25
$ws.Range("J132").Value = 2835.2307
$ws.Range("K132").Value = 1574907.75
$ws.Range("L132").Value = 8505.6921
$ws.Range("M132").Value = -1572377.75
$ws.Range("N132").Value = -13565.6921

$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display
$ws.Range("H96").Value = 1942.5555
$ws.Range("I96").Value = 1863.1666
$ws.Range("K96").Value = 1863.1666
$ws.Range("M96").Value = -490.1666

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 1533.1613
$ws.Range("I126").Value = 1116.4615
$ws.Range("K126").Value = 3349.3845
$ws.Range("M126").Value = -879.3844999999997

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 1205.4286
$ws.Range("I132").Value = 757.4194
$ws.Range("J132").Value = 2468
$ws.Range("K132").Value = 2272.2582
$ws.Range("L132").Value = 7404
$ws.Range("M132").Value = 257.7417999999998
$ws.Range("N132").Value = -12464

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 22941626
$ws.Range("I136").Value = 29494178
$ws.Range("J136").Value = 7691.5
$ws.Range("K136").Value = 88482534
$ws.Range("L136").Value = 23074.5
$ws.Range("M136").Value = -88479984
$ws.Range("N136").Value = -28174.5
